# Commit: "add program xxptrp07.p."
#
# The workbook's physical xl/worksheets/sheet1.xml maps (via
# xl/_rels/workbook.xml.rels: rId1 -> sheet1.xml) to the worksheet named
# "gplbldmt.p" (the first, active tab) -- NOT the sheet literally named
# "Sheet1". We target it by name to be unambiguous.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gplbldmt.p")
$ws.Activate()

# New field/program/label rows documenting program xxptrp07.p, appended
# after the existing data (which currently ends at row 29).
$newRows = @(
    @("in_site",        "xxptrp07.p", "SITE"),
    @("t_lddet_loc",     "xxptrp07.p", "LOCATION"),
    @("t_lddet_part",    "xxptrp07.p", "ITEM_NUMBER"),
    @("pt_desc1",        "xxptrp07.p", "DESCRIPTION"),
    @("pt_draw",         "xxptrp07.p", "DRAWING"),
    @("t_sct_abc",       "xxptrp07.p", "ABC"),
    @("t_lddet_qty",     "xxptrp07.p", "QUANTITY_ON_HAND"),
    @("t_sct_um",        "xxptrp07.p", "UNITS"),
    @("t_sct_std_as_of", "xxptrp07.p", "GL_COST"),
    @("ext_std",         "xxptrp07.p", "GL_COST_TOTAL"),
    @("t_part_type",     "xxptrp07.p", "TYPE"),
    @("t_acct",          "xxptrp07.p", "ACCTS"),
    @("t_sub",           "xxptrp07.p", "SUB-ACCOUNT")
)

$startRow = 30
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
}

# Move the view: the new content scrolls the frozen window down so row 17
# is the first visible row below the frozen header, with A35 selected.
$ws.Range("A35").Select()
